$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - copy formatting from the neighboring
# header cell (G1) so it picks up the same bold/border/alignment style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill the new column's data rows with 0, matching the existing rows.
$ws.Range("H2:H5").Value = 0
